# Add 21th_day_test and 60th_day_test into Estrus.
# Insert a new header row at the top of the "基本資料" sheet, shifting all
# existing data down by one row, then populate the header row with the
# English column names and update a few derived "{field}" placeholder
# messages so they reflect the actual (renamed) field names.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("基本資料")

# Shift all existing rows down by one to make room for the new header row.
$ws.Rows(1).Insert()

# Populate the new header row.
$ws.Range("A1").Value = "Breed"
$ws.Range("B1").Value = "ID"
$ws.Range("C1").Value = "confusing_note"
$ws.Range("D1").Value = "Birthday"
$ws.Range("E1").Value = "Sire"
$ws.Range("F1").Value = "Dam"
$ws.Range("G1").Value = "reg_id"
$ws.Range("H1").Value = "Chinese_name"
$ws.Range("I1").Value = "Gender"
$ws.Range("J1").Value = "註釋"

# Rows below shifted down by one; update the validation-message cells that
# referenced the generic "{field}" placeholder so they show the real
# (English) field names that replaced the old header labels.
$ws.Range("J20").Value = "['Sire 不能為空值', 'Dam 不能為空值', '父畜品種不在常見名單內', '母畜品種不在常見名單內']"
$ws.Range("J28").Value = "['Chinese_name 不能為空值', 'Sire 不能為空值', 'Dam 不能為空值', '父畜品種不在常見名單內', '母畜品種不在常見名單內', '不允許有相近耳號']"
$ws.Range("J66").Value = "['Birthday 不能為空值', '需要有子代的生日才能設定親代', '需要有子代的生日才能設定親代', '不允許有相近耳號']"
